$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("samples_retained")

# New row 15: enterface_db dataset details
$ws.Range("A15").Value = "enterface_db"
$ws.Range("B15").Value = "spon."
$ws.Range("C15").Value = 213
$ws.Range("D15").Value = 1080
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = "English"
$ws.Range("G15").Formula = '=IF(OR(ISBLANK(C15), ISBLANK(D15),ISBLANK(E15)), "", SUM(C15:E15))'
$ws.Range("H15").Value = "elicited in lab; "

# Row 16: carry the running total formula one row further (empty data row)
$ws.Range("G16").Formula = '=IF(OR(ISBLANK(C16), ISBLANK(D16),ISBLANK(E16)), "", SUM(C16:E16))'

# Update selection to match final state
$ws.Range("G16").Select()
